# Auto-generated edit script: update crypto price/volume table (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.825.20'
$ws.Range('E2').Value = '  +6.61%  '
$ws.Range('D3').Value = '1.740.16'
$ws.Range('E3').Value = '  +5.33%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'227.74"
$ws.Range('E5').Value = '  +4.23%  '
$ws.Range('D6').Value = "'0.5463"
$ws.Range('E6').Value = '  +3.95%  '
$ws.Range('D7').Value = "'1.004"
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = "'0.2779"
$ws.Range('E8').Value = '  +3.97%  '
$ws.Range('D9').Value = "'0.06755"
$ws.Range('E9').Value = '  +6.08%  '
$ws.Range('D10').Value = "'22.01"
$ws.Range('E10').Value = '  +7.10%  '
$ws.Range('D11').Value = "'0.07781"
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('D12').Value = "'4.705"
$ws.Range('E12').Value = '  +2.40%  '
$ws.Range('D13').Value = '1.753.50'
$ws.Range('E13').Value = '  +3.65%  '
$ws.Range('D14').Value = '1.980.20'
$ws.Range('E14').Value = '  +5.32%  '
$ws.Range('D15').Value = "'0.5986"
$ws.Range('E15').Value = '  +6.81%  '
$ws.Range('D16').Value = '0.0₅8418'
$ws.Range('E16').Value = '  +2.29%  '
$ws.Range('D17').Value = "'69.09"
$ws.Range('E17').Value = '  +5.57%  '
$ws.Range('D18').Value = '27.829.14'
$ws.Range('E18').Value = '  +6.59%  '
$ws.Range('D19').Value = "'225.63"
$ws.Range('E19').Value = '  +18.03%  '
$ws.Range('D20').Value = "'4.848"
$ws.Range('E20').Value = '  +3.38%  '
$ws.Range('D21').Value = "'1.004"
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('E22').Value = '  +5.50%  '
$ws.Range('D23').Value = "'6.247"
$ws.Range('E23').Value = '  +4.60%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').Value = "'146.55"
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').Value = "'0.1252"
$ws.Range('E26').Value = '  +4.20%  '
$ws.Range('D27').Value = "'1.677"
$ws.Range('E27').Value = '  +12.02%  '
$ws.Range('D28').Value = "'7.471"
$ws.Range('E28').Value = '  +3.02%  '
$ws.Range('D29').Value = "'17.21"
$ws.Range('E29').Value = '  +7.86%  '
$ws.Range('D30').Value = "'0.05676"
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('D31').Value = "'1.314"
$ws.Range('E31').Value = '  +3.27%  '
$ws.Range('D32').Value = "'3.700"
$ws.Range('E32').Value = '  +5.80%  '
$ws.Range('D33').Value = "'3.529"
$ws.Range('E33').Value = '  +4.43%  '
$ws.Range('D34').Value = "'1.690"
$ws.Range('E34').Value = '  +7.05%  '
$ws.Range('D35').Value = "'0.9782"
$ws.Range('E35').Value = '  +3.48%  '
$ws.Range('D36').Value = "'2.858"
$ws.Range('E36').Value = '  +2.03%  '
$ws.Range('E37').Value = '  +1.78%  '
$ws.Range('D38').Value = "'0.5980"
$ws.Range('E38').Value = '  +3.53%  '
$ws.Range('D39').Value = "'0.01670"
$ws.Range('E39').Value = '  +4.84%  '
$ws.Range('D40').Value = "'5.992"
$ws.Range('E40').Value = '  +0.34%  '
$ws.Range('D41').Value = "'0.8513"
$ws.Range('E41').Value = '  +1.27%  '
$ws.Range('D42').Value = '1.047.59'
$ws.Range('E42').Value = '  +2.32%  '
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('D44').Value = "'102.13"
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('D45').Value = '1.886.52'
$ws.Range('E45').Value = '  +5.28%  '
$ws.Range('D46').Value = "'59.56"
$ws.Range('E46').Value = '  +1.69%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'8.303"
$ws.Range('E47').Value = '  +3.46%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = "'0.4441"
$ws.Range('E48').Value = '  +2.25%  '
$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').Value = "'1.008"
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.05321"
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = "'0.1013"
$ws.Range('E51').Value = '  +3.83%  '

Write-Output "Applied 105 cell updates (37 forced text)"
